# fix: quản lý đơn nghỉ phép, xuất báo cáo
#
# Updates the export-leave-template.xlsx:
#  - Rename placeholder tokens in the repeating-row template (D6, E6, I6)
#    so the exporter uses the "converted"/named variants.
#  - Widen column E (width 19.33 -> 22.56) to fit the longer placeholder.
#  - Update the sheet's saved view state (scrolled/selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string placeholders used by the forEach template row ---
$ws.Range("D6").Value = '${item.startDayConvert}'
$ws.Range("E6").Value = '${item.endDayConvert}'
$ws.Range("I6").Value = '${item.isActiveName}'

# --- Widen column E to fit the new placeholder text ---
# (target raw width 22.5555555555556; COM ColumnWidth is quantized to 1/6-
# character steps by this host, so 21.65 is the nearest input that lands on
# the closest achievable stored width, 22.5)
$ws.Columns.Item(5).ColumnWidth = 21.65

# --- Update sheet view: scroll position + active selection ---
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I9").Select()
